$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: D7 <-> E7 swap
$ws.Range("D7").Value = "-"
$ws.Range("E7").Value = "['MEC-3B-Tec. Mat. Não Metal.', 'MEC-3B-Tec. Mat. Não Metal.']"

# Row 18: C18 gets new value, D18 becomes "-"
$ws.Range("C18").Value = "[-, 'MEC-1NB-Tec. Mat. Não Metal.', -, -]"
$ws.Range("D18").Value = "-"

# Row 19: C19 gets new value, D19 becomes "-"
$ws.Range("C19").Value = "[-, 'MEC-1NB-Tec. Mat. Não Metal.', -, -]"
$ws.Range("D19").Value = "-"

# Row 21: B21 gets new value, D21 becomes "-"
$ws.Range("B21").Value = "[-, 'MEC-1NB-Tec. Mat. Não Metal.', -, -]"
$ws.Range("D21").Value = "-"
